{"js": "// Update the date line (first paragraph) from \"2025-01-02 Thursday\" to\n// \"2025-01-03 Friday\", then replace every arithmetic-fact string inside the\n// 20x5 practice table with its new value (same cell positions, new numbers).\n\nconst body = context.document.body;\n\n// --- 1. Date paragraph -------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst datePara = paragraphs.items[0];\ndatePara.load(\"text\");\nawait context.sync();\n\nif (datePara.text.indexOf(\"2025-01-02 Thursday\") !== -1) {\n  datePara.insertText(\"2025-01-03 Friday\", Word.InsertLocation.replace);\n} else {\n  // Fallback: if the exact phrase isn't the whole paragraph text, do a\n  // scoped search-and-replace instead so we don't disturb anything else.\n  const results = datePara.search(\"2025-01-02 Thursday\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"2025-01-03 Friday\", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\n// --- 2. Practice table ---------------------------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// New values, row-major, matching the existing 20 rows x 5 columns layout.\nconst newValues = [\n  [\"22+55=77\", \"6+3=9\", \"79-34=45\", \"30+58=88\", \"16+66=82\"],\n  [\"68-14=54\", \"84-28=56\", \"80-40=40\", \"95-88=7\", \"63-20=43\"],\n  [\"86-62=24\", \"83+9=92\", \"48-47=1\", \"56-23=33\", \"27+9=36\"],\n  [\"35-17=18\", \"22+27=49\", \"98-71=27\", \"20+4=24\", \"72-11=61\"],\n  [\"45-26=19\", \"92-57=35\", \"12+57=69\", \"66-62=4\", \"7+66=73\"],\n  [\"19+75=94\", \"51+18=69\", \"30+55=85\", \"77-3=74\", \"82+11=93\"],\n  [\"3+69=72\", \"6+63=69\", \"86-67=19\", \"81-44=37\", \"69-19=50\"],\n  [\"59-54=5\", \"16+49=65\", \"29+14=43\", \"17+17=34\", \"76-63=13\"],\n  [\"21+10=31\", \"27-15=12\", \"52+21=73\", \"28-2=26\", \"1+94=95\"],\n  [\"60-59=1\", \"49-33=16\", \"9+0=9\", \"20+54=74\", \"34+26=60\"],\n  [\"94+5=99\", \"77-3=74\", \"4+39=43\", \"13+37=50\", \"78+10=88\"],\n  [\"60-49=11\", \"47-43=4\", \"26+55=81\", \"85-4=81\", \"0+61=61\"],\n  [\"91-71=20\", \"97-20=77\", \"66-41=25\", \"38-29=9\", \"82-17=65\"],\n  [\"85-7=78\", \"49+34=83\", \"40-15=25\", \"94-44=50\", \"20+29=49\"],\n  [\"77-45=32\", \"49-6=43\", \"81+8=89\", \"25-3=22\", \"45+51=96\"],\n  [\"78-41=37\", \"96-90=6\", \"38-3=35\", \"23+0=23\", \"89-73=16\"],\n  [\"0+26=26\", \"85-31=54\", \"67-16=51\", \"61-48=13\", \"32+49=81\"],\n  [\"80-47=33\", \"12+62=74\", \"72-67=5\", \"16+72=88\", \"55+40=95\"],\n  [\"56+0=56\", \"73-10=63\", \"86-77=9\", \"7+10=17\", \"9+5=14\"],\n  [\"51-48=3\", \"46+45=91\", \"80-19=61\", \"88-14=74\", \"95-80=15\"],\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the date line, then replace every arithmetic-fact string inside\n# the 20x5 practice table with its new value (same cell positions, new\n# numbers).\n\n$d = $word.ActiveDocument\n\n# --- 1. Date paragraph ---------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\n    \"2025-01-02 Thursday\",  # FindText\n    $false,        # MatchCase\n    $false,        # MatchWholeWord\n    $false,        # MatchWildcards\n    $false,        # MatchSoundsLike\n    $false,        # MatchAllWordForms\n    $true,         # Forward\n    1,             # Wrap (wdFindContinue)\n    $false,        # Format\n    \"2025-01-03 Friday\",  # ReplaceWith\n    2              # Replace (wdReplaceAll)\n) | Out-Null\n\n# --- 2. Practice table -----------------------------------------------------\n$table = $d.Tables.Item(1)\n\n# New values, row-major, matching the existing 20 rows x 5 columns layout.\n$newValues = @(\n    @(\"22+55=77\", \"6+3=9\", \"79-34=45\", \"30+58=88\", \"16+66=82\"),\n    @(\"68-14=54\", \"84-28=56\", \"80-40=40\", \"95-88=7\", \"63-20=43\"),\n    @(\"86-62=24\", \"83+9=92\", \"48-47=1\", \"56-23=33\", \"27+9=36\"),\n    @(\"35-17=18\", \"22+27=49\", \"98-71=27\", \"20+4=24\", \"72-11=61\"),\n    @(\"45-26=19\", \"92-57=35\", \"12+57=69\", \"66-62=4\", \"7+66=73\"),\n    @(\"19+75=94\", \"51+18=69\", \"30+55=85\", \"77-3=74\", \"82+11=93\"),\n    @(\"3+69=72\", \"6+63=69\", \"86-67=19\", \"81-44=37\", \"69-19=50\"),\n    @(\"59-54=5\", \"16+49=65\", \"29+14=43\", \"17+17=34\", \"76-63=13\"),\n    @(\"21+10=31\", \"27-15=12\", \"52+21=73\", \"28-2=26\", \"1+94=95\"),\n    @(\"60-59=1\", \"49-33=16\", \"9+0=9\", \"20+54=74\", \"34+26=60\"),\n    @(\"94+5=99\", \"77-3=74\", \"4+39=43\", \"13+37=50\", \"78+10=88\"),\n    @(\"60-49=11\", \"47-43=4\", \"26+55=81\", \"85-4=81\", \"0+61=61\"),\n    @(\"91-71=20\", \"97-20=77\", \"66-41=25\", \"38-29=9\", \"82-17=65\"),\n    @(\"85-7=78\", \"49+34=83\", \"40-15=25\", \"94-44=50\", \"20+29=49\"),\n    @(\"77-45=32\", \"49-6=43\", \"81+8=89\", \"25-3=22\", \"45+51=96\"),\n    @(\"78-41=37\", \"96-90=6\", \"38-3=35\", \"23+0=23\", \"89-73=16\"),\n    @(\"0+26=26\", \"85-31=54\", \"67-16=51\", \"61-48=13\", \"32+49=81\"),\n    @(\"80-47=33\", \"12+62=74\", \"72-67=5\", \"16+72=88\", \"55+40=95\"),\n    @(\"56+0=56\", \"73-10=63\", \"86-77=9\", \"7+10=17\", \"9+5=14\"),\n    @(\"51-48=3\", \"46+45=91\", \"80-19=61\", \"88-14=74\", \"95-80=15\")\n)\n\nfor ($r = 1; $r -le $newValues.Count; $r++) {\n    $row = $newValues[$r - 1]\n    for ($c = 1; $c -le $row.Count; $c++) {\n        $table.Cell($r, $c).Range.Text = $row[$c - 1]\n    }\n}\n"}
